$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet")

$names = @(
    "Bitcoin",
    "Ethereum",
    "Tether",
    "BNB",
    "XRP",
    "Solana",
    "Dogecoin",
    "Cardano",
    "USDC",
    "Lido Staked ETH",
    "TRON",
    "Avalanche",
    "Shiba Inu",
    "Stellar",
    "Chainlink",
    "Polkadot",
    "Wrapped Bitcoin",
    "Bitcoin Cash",
    "Litecoin",
    "Uniswap",
    "NEAR Protocol",
    "UNUS SED LEO",
    "Internet Computer",
    "VeChain",
    "Multi Collateral DAI",
    "Ethereum Classic",
    "Crypto.com Coin",
    "Sp8de",
    "Hedera Hashgraph",
    "Artificial Superintelligence Alliance",
    "Filecoin",
    "Algorand",
    "Stacks",
    "Monero",
    "OKB",
    "Aave",
    "Fantom",
    "The Graph",
    "THETA",
    "Injective",
    "Cosmos",
    "THORChain",
    "EOS",
    "Maker",
    "IOTA",
    "MANTRA DAO",
    "The Sandbox",
    "Gala",
    "Lido DAO",
    "Flow",
    "Arweave",
    "Quant",
    "Tezos",
    "Neo",
    "Polygon",
    "KuCoin Token",
    "Helium",
    "Bitcoin SV",
    "Raydium",
    "Axie Infinity",
    "Decentraland",
    "MultiversX",
    "AIOZ Network",
    "Conflux",
    "Zcash",
    "GateToken",
    "Reserve Rights",
    "Chiliz",
    "Mina",
    "Curve DAO Token",
    "eCash",
    "Akash Network",
    "Pendle",
    "PancakeSwap",
    "XinFin Network",
    "Oasis",
    "Nexo",
    "Amp",
    "FTX Token",
    "SuperVerse",
    "Nervos Network",
    "Compound",
    "Kava",
    "Gnosis",
    "Kusama",
    "Dash",
    "1inch Network",
    "Synthetix",
    "Zilliqa",
    "Holo",
    "Bitcoin Gold",
    "WOO",
    "Telcoin",
    "Enjin Coin",
    "Theta Fuel",
    "Livepeer",
    "IoTeX",
    "DeXe",
    "JUST",
    "Celo"
)

$prices = @(
    95706.12940484812,
    3647.614127816061,
    1.00029194544417,
    751.5138343936801,
    2.517656103469206,
    236.1708739410761,
    0.4124796933033013,
    1.196424070616533,
    0.999887564710715,
    3620.467876837743,
    0.3745805180736099,
    51.45218344200527,
    0.0000299999301083,
    0.5061716247575836,
    23.74415752206779,
    9.753780387062024,
    95444.4697948156,
    564.7635566181027,
    128.4899641654699,
    14.72264708288165,
    7.230190371054206,
    9.378120456024803,
    14.30680728225584,
    0.0713035356852022,
    0.9997725402354023,
    35.5635463814456,
    0.1989262583227424,
    0.6305730521682377,
    0.3230489095306408,
    1.899088648975919,
    7.50177219157705,
    0.5225528938049699,
    2.568231338139122,
    199.1759878162702,
    60.76082907329035,
    238.8627715284042,
    1.234611448289417,
    0.3119893596470023,
    2.92134492588411,
    32.79218321236868,
    9.50312939339447,
    6.813167763105839,
    1.439464472980603,
    2320.27761085187,
    0.5799393202739836,
    3.865202862651067,
    0.7688972455961697,
    0.0511291384015356,
    2.040327398783852,
    1.153409346194165,
    27.18831590262465,
    142.5754269695325,
    1.689243969923392,
    23.09408529776423,
    0.7042263434138653,
    12.86631767220504,
    8.880404548421943,
    76.08217931632423,
    5.070599069983644,
    9.201375304239761,
    0.6905579344804675,
    52.56985537178777,
    1.111713337277815,
    0.2556315291398551,
    72.41205538679591,
    12.70939935618371,
    0.0209272965523953,
    0.1193086562992418,
    0.9173660555248648,
    0.8747312753173011,
    0.0000539320011504,
    4.256929732800584,
    6.36759257895398,
    3.52424150719294,
    0.07811543886775341,
    0.1358193749089508,
    1.479580060794976,
    0.0117482380443289,
    2.724066151435012,
    1.703215909235181,
    0.0179734077664143,
    87.33242034194339,
    0.6889598030884678,
    277.4327337683155,
    43.81442483102403,
    56.78322048768199,
    0.5152824568864862,
    2.86570576768944,
    0.0342119069593094,
    0.0036875115419074,
    35.54908280363081,
    0.3353212877544386,
    0.0067318032822877,
    0.3460187290980695,
    0.0902310787773416,
    16.42672815791765,
    0.0631239193988599,
    10.3290787488951,
    0.0583069147672488,
    1.034056470904978
)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
    $ws.Cells.Item($row, 3).Value = $prices[$i]
}

$wsStats = $wb.Worksheets.Item("Statistics")
$wsStats.Range("A2").Value = 2041.7
$wsStats.Range("B2").Value = 2.5
$wsStats.Range("C2").Value = 13441
